$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("21:21").Insert()

$ws.Range("A21").Value = 3
$ws.Range("B21").Value = 'Femacal de La Calera'
$ws.Range("C21").Value = 'Coquimbo'
$ws.Range("D21").Value = 44473
$ws.Range("E21").Value = 5
$ws.Range("F21").Value = 'Fruta'
$ws.Range("G21").Value = 100101
$ws.Range("H21").Value = 'Berries'
$ws.Range("I21").Value = 100101001
$ws.Range("J21").Value = 'Arándano (blue)'
$ws.Range("K21").Value = 'Sin especificar'
$ws.Range("L21").Value = 'Primera'
$ws.Range("M21").Value = 75
$ws.Range("N21").Value = 11000
$ws.Range("O21").Value = 11000
$ws.Range("P21").Value = 11000
$ws.Range("Q21").Value = '$/bandeja 12 canastillos 125 gramos'
$ws.Range("R21").Value = 'Provincia de Quillota'
$ws.Range("S21").Value = 7333
$ws.Range("T21").Value = 1.5
